# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet named "2022-Q4" right after "总计" (i.e. right
#    before the existing "2022-Q1" sheet), cloning the layout/formatting of the
#    "2022-Q1" sheet (which holds the same 8-column fund-holding table shape),
#    then overwrite its data rows with the 2022-Q4 fund holdings (5 funds).
# 2. Add a new "2022-Q4" summary row to the "总计" totals sheet, right after the
#    header row, and keep the rest of the historical rows intact.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" worksheet
# ---------------------------------------------------------------------------

$q1sheet = $wb.Worksheets.Item("2022-Q1")

# Worksheet.Copy(Before) clones the sheet (formatting + content) and places
# the clone immediately before $q1sheet, i.e. right after "总计".
$q1sheet.Copy($q1sheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Clear out the old 2022-Q1 data row that got copied along with the sheet, and
# write the five 2022-Q4 fund rows in its place. Rows 3-6 are brand new, so
# pick up the row-2 (data row) formatting for column A explicitly afterwards.
$newSheet.Rows.Item(2).ClearContents()

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "'016285"
$newSheet.Cells.Item(2, 3).Value = "汇丰晋信龙头优势混合A"
$newSheet.Cells.Item(2, 4).Value = "'1.61"
$newSheet.Cells.Item(2, 5).Value = "'58.20"
$newSheet.Cells.Item(2, 6).Value = "'3.38"
$newSheet.Cells.Item(2, 7).Value = "'0.0544"
$newSheet.Cells.Item(2, 8).Value = 1

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "'016174"
$newSheet.Cells.Item(3, 3).Value = "汇丰晋信策略优选混合A"
$newSheet.Cells.Item(3, 4).Value = "'1.78"
$newSheet.Cells.Item(3, 5).Value = "'74.92"
$newSheet.Cells.Item(3, 6).Value = "'2.57"
$newSheet.Cells.Item(3, 7).Value = "'0.0457"
$newSheet.Cells.Item(3, 8).Value = 3

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "'016175"
$newSheet.Cells.Item(4, 3).Value = "汇丰晋信策略优选混合C"
$newSheet.Cells.Item(4, 4).Value = "'0.40"
$newSheet.Cells.Item(4, 5).Value = "'74.92"
$newSheet.Cells.Item(4, 6).Value = "'2.57"
$newSheet.Cells.Item(4, 7).Value = "'0.0103"
$newSheet.Cells.Item(4, 8).Value = 3

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "'016286"
$newSheet.Cells.Item(5, 3).Value = "汇丰晋信龙头优势混合C"
$newSheet.Cells.Item(5, 4).Value = "'0.26"
$newSheet.Cells.Item(5, 5).Value = "'58.20"
$newSheet.Cells.Item(5, 6).Value = "'3.38"
$newSheet.Cells.Item(5, 7).Value = "'0.0088"
$newSheet.Cells.Item(5, 8).Value = 1

$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = "'562530"
$newSheet.Cells.Item(6, 3).Value = "华夏中证智选1000价值稳健策略ETF"
$newSheet.Cells.Item(6, 4).Value = "'0.36"
$newSheet.Cells.Item(6, 5).Value = "'96.22"
$newSheet.Cells.Item(6, 6).Value = "'0.98"
$newSheet.Cells.Item(6, 7).Value = "'0.0035"
$newSheet.Cells.Item(6, 8).Value = 2

# Column A keeps the same centered/bold/bordered look on every data row - copy
# that formatting down from row 2 onto the newly-appended rows 3-6.
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: add the 2022-Q4 row to the "总计" summary sheet
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Inserting a row copies the formatting of the row above (the bold header), so
# reset the plain data cells and re-apply column A's normal data-row look.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.12

# Re-number the 0-based index column (A) for the rows that followed.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(7, 1).Value = 5

# Restore the originally-active "2020-Q4" tab (creating/copying sheets moves
# the selection around as a side effect).
$wb.Worksheets.Item("2020-Q4").Activate()
